# "Updated Test Data and imports"
# - CreateProject!B2:B5 data values get shuffled around.
# - CreateProject becomes the active/selected sheet (was UpdateProject).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CreateProject")

$ws1.Range("B2").Value = "data3"
$ws1.Range("B3").Value = "data5"
$ws1.Range("B4").Value = "data2"
$ws1.Range("B5").Value = "data5"

# Make CreateProject the active sheet / selection (UpdateProject loses it).
$null = $ws1.Select()
$null = $ws1.Range("A1").Select()
